# Apply the edit described by the diff:
# 1. Rotate rows 184-187 (row 184's data moves down to row 187, rows 185-187 shift up by one).
# 2. Swap rows 285-286.
# 3. Update every timestamp in column O (rows 2-333) to the new crawl time.
#
# Columns D (ratingAmount) and E (ratingValue) hold real numbers; every other
# column holds text (ids, prices, urls, etc. are all stored as strings even
# when they look numeric). When copying a row via a Value2 array and writing
# it back, Excel-COM will silently reinterpret numeric-looking strings as
# numbers, so those columns are protected with a leading apostrophe (which
# forces text and is stripped back out) before being written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Protect-TextColumns($rowArr) {
    for ($c = 1; $c -le 15; $c++) {
        if ($c -eq 4 -or $c -eq 5) { continue }   # D, E are numeric columns
        $v = $rowArr[1, $c]
        if ($v -ne $null) {
            $rowArr[1, $c] = "'" + $v
        }
    }
    return $rowArr
}

function Get-RowValues($rowNum) {
    return $ws.Range("A$rowNum`:O$rowNum").Value2
}

function Set-RowValues($rowNum, $values) {
    $ws.Range("A$rowNum`:O$rowNum").Value = (Protect-TextColumns $values)
}

# --- rotate rows 184-187 ----------------------------------------------------
$row184 = Get-RowValues 184
$row185 = Get-RowValues 185
$row186 = Get-RowValues 186
$row187 = Get-RowValues 187

Set-RowValues 184 $row185
Set-RowValues 185 $row186
Set-RowValues 186 $row187
Set-RowValues 187 $row184

# --- swap rows 285-286 ------------------------------------------------------
$row285 = Get-RowValues 285
$row286 = Get-RowValues 286

Set-RowValues 285 $row286
Set-RowValues 286 $row285

# --- refresh the scrape timestamp on every data row -------------------------
$ws.Range("O2:O333").Value = "2023-02-05 20:49:26"
